$d = $word.ActiveDocument

function Replace-Merged([string]$searchText, [string]$finalText) {
    # Step 1: find the target text and stomp it with a unique placeholder.
    # This forces the whole matched range (which may span several runs,
    # each separated by <w:proofErr/> spell-check markers) to collapse
    # into a single freshly generated run, dropping the now-orphaned
    # <w:proofErr/> markers along with the old run boundaries.
    $placeholder = "__REPL_" + [guid]::NewGuid().ToString("N") + "__"

    $rng1 = $d.Content
    $found1 = $rng1.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $placeholder, 2)
    if (-not $found1) {
        throw "Could not find text: $searchText"
    }

    # Step 2: swap the placeholder back out for the real final text. Using
    # a two-step replace (instead of replacing $searchText directly with
    # $finalText) avoids a same-text no-op when $finalText happens to
    # equal $searchText (e.g. "${ilos}" -> "${ilos}", only cleaning up
    # the run/proofErr structure).
    $rng2 = $d.Content
    $found2 = $rng2.Find.Execute($placeholder, $false, $false, $false, $false, $false, $true, 1, $false, $finalText, 2)
    if (-not $found2) {
        throw "Could not find placeholder: $placeholder"
    }
}

# departmentFull -> support (kept as its own run, flanked by the
# pre-existing "${" and "}" runs, same as before the edit).
Replace-Merged '${departmentFull}' '${support}'

# Re-split the merged "${support}" run back into three runs -- "${",
# "support", "}" -- mirroring the original run layout, by nudging a
# character-formatting property on just the middle word and reverting it.
# This forces the run to be carved out of its neighbors without picking
# up any <w:proofErr/> wrapping.
$rng = $d.Content
$found = $rng.Find.Execute('${support}', $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw 'Could not find ${support} placeholder span'
}
$s = $rng.Start
$e = $rng.End
$mid = $d.Range($s + 2, $e - 1)
$mid.Bold = 1
$mid.Bold = 0

# The remaining ${...} placeholders just need their three runs (prefix,
# name, suffix) merged back into one clean run with the same text, which
# also clears out the <w:proofErr/> wrapping around the name.
Replace-Merged '${ilos}' '${ilos}'
Replace-Merged '${budgetSource}' '${budgetSource}'
Replace-Merged '${sig_cscp}' '${sig_cscp}'
Replace-Merged '${sig_csca}' '${sig_csca}'
Replace-Merged '${sig_sscp}' '${sig_sscp}'
Replace-Merged '${sig_dean}' '${sig_dean}'

Write-Output "done"
